$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Release"

$releaseValues = @(1, 1, 2, 1, 1, 2, 1, 2, 2, 1)
for ($i = 0; $i -lt $releaseValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $releaseValues[$i]
}

$ws.Range("D12").Select()
